$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.158.44'
$ws.Range('E2').Value = '  -1.53%  '

$ws.Range('D3').Value = '2.178.55'
$ws.Range('E3').Value = '  -1.76%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = "'250.10"
$ws.Range('E5').Value = '  -0.90%  '

$ws.Range('D6').Value = "'0.612"

$ws.Range('D7').Value = "'66.37"
$ws.Range('E7').Value = '  -6.75%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').Value = "'0.591"
$ws.Range('E9').Value = '  -0.60%  '

$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').Value = "'37.15"
$ws.Range('E10').Value = '  -8.61%  '

$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = "'59.11"
$ws.Range('E11').Value = '  +1.53%  '

$ws.Range('D12').Value = "'0.0935"
$ws.Range('E12').Value = '  -3.23%  '

$ws.Range('E13').Value = '  -0.79%  '

$ws.Range('D14').Value = "'6.87"
$ws.Range('E14').Value = '  -4.82%  '

$ws.Range('D15').Value = '2.505.22'
$ws.Range('E15').Value = '  -1.64%  '

$ws.Range('D16').Value = "'14.28"
$ws.Range('E16').Value = '  -4.53%  '

$ws.Range('D17').Value = "'0.847"
$ws.Range('E17').Value = '  -3.35%  '

$ws.Range('D18').Value = '2.165.97'
$ws.Range('E18').Value = '  -2.34%  '

$ws.Range('D19').Value = '41.141.53'
$ws.Range('E19').Value = '  -1.57%  '

$ws.Range('D20').Value = '0.0₃0945'
$ws.Range('E20').Value = '  -1.72%  '

$ws.Range('D21').Value = "'71.60"
$ws.Range('E21').Value = '  -1.40%  '

$ws.Range('D22').Value = "'6.05"
$ws.Range('E22').Value = '  -2.64%  '

$ws.Range('D23').Value = "'230.03"
$ws.Range('E23').Value = '  -2.08%  '

$ws.Range('D24').Value = "'2.01"
$ws.Range('E24').Value = '  -2.41%  '

$ws.Range('E25').Value = '  -6.65%  '

$ws.Range('E26').Value = '  -0.03%  '

$ws.Range('D27').Value = "'11.18"
$ws.Range('E27').Value = '  +0.51%  '

$ws.Range('D28').Value = "'2.40"
$ws.Range('E28').Value = '  -5.66%  '

$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').Value = "'3.68"
$ws.Range('E29').Value = '  -4.17%  '

$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = "'166.90"
$ws.Range('E30').Value = '  -1.87%  '

$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = "'2.04"
$ws.Range('E31').Value = '  -7.51%  '

$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'20.13"
$ws.Range('E32').Value = '  -3.01%  '

$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = "'0.120"
$ws.Range('E33').Value = '  -1.88%  '

$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = "'5.71"
$ws.Range('E34').Value = '  +3.32%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = "'0.0750"
$ws.Range('E35').Value = '  +1.39%  '

$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = "'0.121"
$ws.Range('E36').Value = '  -1.03%  '

$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').Value = "'4.51"
$ws.Range('E37').Value = '  -3.58%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = "'4.06"
$ws.Range('E38').Value = '  +1.10%  '

$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = "'24.99"
$ws.Range('E39').Value = '  -5.76%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.0306"
$ws.Range('E40').Value = '  -0.21%  '

$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D41').Value = "'5.45"
$ws.Range('E41').Value = '  +15.73%  '

$ws.Range('B42').Value = 'LidoDAOToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D42').Value = "'2.21"
$ws.Range('E42').Value = '  -2.75%  '

$ws.Range('B43').Value = 'THORChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D43').Value = "'5.53"
$ws.Range('E43').Value = '  -6.79%  '

$ws.Range('D44').Value = "'11.40"
$ws.Range('E44').Value = '  -10.33%  '

$ws.Range('B45').Value = 'MultiversX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D45').Value = "'60.47"
$ws.Range('E45').Value = '  -5.89%  '

$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = "'8.51"
$ws.Range('E46').Value = '  -2.50%  '

$ws.Range('D47').Value = "'0.190"
$ws.Range('E47').Value = '  -6.40%  '

$ws.Range('B48').Value = 'BinanceUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D48').Value = "'1.00"
$ws.Range('E48').Value = '  -0.07%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = "'0.0987"
$ws.Range('E49').Value = '  -3.23%  '

$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').Value = "'1.14"
$ws.Range('E50').Value = '  -2.53%  '

$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').Value = "'4.39"
$ws.Range('E51').Value = '  -6.07%  '
